# Update countries & provincias Spain
# Applies updated COVID-19 stats to the "Pais" sheet, which causes a
# handful of countries (Honduras, Bulgaria, Guatemala, Hong Kong, Tunez,
# Mayotte, Kirguistan, Consejo Danes para los Refugiados) to swap/shift
# position because the sheet is kept sorted by total cases (column B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 42: Corea del Sur ---
$ws.Range("B42").Value = 10909
$ws.Range("C42").Value = 35
$ws.Range("D42").Value = 9632
$ws.Range("E42").Value = 1021

# --- Row 53: Australia ---
$ws.Range("B53").Value = 6948
$ws.Range("C53").Value = 7
$ws.Range("D53").Value = 6167
$ws.Range("E53").Value = 684

# --- Rows 79-80: Honduras overtakes Bulgaria ---
$ws.Range("A79").Value = "Honduras"
$ws.Range("B79").Value = 1972
$ws.Range("C79").Value = 142
$ws.Range("D79").Value = 203
$ws.Range("E79").Value = 1661
$ws.Range("F79").Value = 10
$ws.Range("H79").Value = 108

$ws.Range("A80").Value = "Bulgaria"
$ws.Range("B80").Value = 1965
$ws.Range("C80").Value = 0
$ws.Range("D80").Value = 444
$ws.Range("E80").Value = 1430
$ws.Range("F80").Value = 58
$ws.Range("H80").Value = 91

# --- Row 87: Nueva Zelanda ---
$ws.Range("B87").Value = 1497
$ws.Range("C87").Value = 3
$ws.Range("D87").Value = 1386
$ws.Range("E87").Value = 90

# --- Rows 94-99: Guatemala overtakes Hong Kong, Tunez, Mayotte,
#     Kirguistan and Consejo Danes para los Refugiados (each shifts
#     down one row) ---
$ws.Range("A94").Value = "Guatemala"
$ws.Range("B94").Value = 1052
$ws.Range("C94").Value = 85
$ws.Range("D94").Value = 110
$ws.Range("E94").Value = 916
$ws.Range("F94").Value = 5
$ws.Range("G94").Value = 2
$ws.Range("H94").Value = 26

$ws.Range("A95").Value = "Hong Kong"
$ws.Range("B95").Value = 1048
$ws.Range("C95").Value = 0
$ws.Range("D95").Value = 982
$ws.Range("E95").Value = 62
$ws.Range("F95").Value = 1
$ws.Range("H95").Value = 4

$ws.Range("A96").Value = "Tunez"
$ws.Range("B96").Value = 1032
$ws.Range("C96").Value = 0
$ws.Range("D96").Value = 700
$ws.Range("E96").Value = 287
$ws.Range("F96").Value = 19
$ws.Range("H96").Value = 45

$ws.Range("A97").Value = "Mayotte"
$ws.Range("B97").Value = 1023
$ws.Range("C97").Value = 0
$ws.Range("D97").Value = 492
$ws.Range("E97").Value = 520
$ws.Range("F97").Value = 9
$ws.Range("H97").Value = 11

$ws.Range("A98").Value = "Kirguistan"
$ws.Range("B98").Value = 1002
$ws.Range("C98").Value = 0
$ws.Range("D98").Value = 675
$ws.Range("E98").Value = 315
$ws.Range("F98").Value = 13
$ws.Range("H98").Value = 12

$ws.Range("A99").Value = "Consejo Danes para los Refugiados"
$ws.Range("B99").Value = 991
$ws.Range("C99").Value = 0
$ws.Range("D99").Value = 136
$ws.Range("E99").Value = 814
$ws.Range("F99").Value = 0
$ws.Range("H99").Value = 41

# --- Row 104: Sri Lanka ---
$ws.Range("B104").Value = 863
$ws.Range("C104").Value = 7
$ws.Range("E104").Value = 533
